$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 10: ORM record for Prototype Study Report (PSR) review.
# No (A10) is already 9 from the template; leave it untouched.

# B10: "ORM No" stored as text (keeps the leading zero) with a
# right-aligned, text-formatted style (numFmtId 49 "@").
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").HorizontalAlignment = -4152
$ws.Range("B10").Value = "06051809"

# C10: Review Document
$ws.Range("C10").Value = "Prototype_Draft.doc"

# D10: Author
$ws.Range("D10").Value = "GaoZuYi"

# E10: Issue Date (5/6/2018 -> serial 43226), reusing the date style
# already used by the column (copy from E9 so no new style is created).
$ws.Range("E9").Copy($ws.Range("E10"))
$ws.Range("E10").Value = 43226

# Move the active selection, matching the saved workbook state.
$ws.Range("J17").Select()
